$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary rows 10-12: update score totals and row-label styles ---
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 60
$ws.Range("C12").Value = -12
$ws.Range("E12").Value = "48/112"

# --- Remove third answer block (columns G:H) entirely ---
$ws.Range("G15:H40").Clear()

# --- Remove second answer block (columns D:E) for rows 19-40 (kept for 16-18) ---
$ws.Range("D19:E40").Clear()

# --- Fill in Student Ans for first block (col A) and update col D for rows 16-18 ---
$ws.Range("A16").Value = "Option B"
$ws.Range("A16").Style = "incorrectStyle"
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"
$ws.Range("A17").Value = "Option B"
$ws.Range("A17").Style = "incorrectStyle"
$ws.Range("D17").Value = "Option B"
$ws.Range("D17").Style = "incorrectStyle"
$ws.Range("A18").Value = "Option B"
$ws.Range("A18").Style = "correctStyle"
$ws.Range("D18").Value = "Option B"
$ws.Range("D18").Style = "incorrectStyle"
$ws.Range("A19").Value = "Option C"
$ws.Range("A19").Style = "correctStyle"
$ws.Range("A20").Value = "Option B"
$ws.Range("A20").Style = "correctStyle"
$ws.Range("A21").Value = "Option C"
$ws.Range("A21").Style = "correctStyle"
$ws.Range("A22").Value = "Option D"
$ws.Range("A22").Style = "correctStyle"
$ws.Range("A23").Value = "Option B"
$ws.Range("A23").Style = "incorrectStyle"
$ws.Range("A24").Value = "Option C"
$ws.Range("A24").Style = "incorrectStyle"
$ws.Range("A25").Value = "Option A"
$ws.Range("A25").Style = "correctStyle"
$ws.Range("A26").Value = "Option B"
$ws.Range("A26").Style = "incorrectStyle"
$ws.Range("A27").Value = "Option A"
$ws.Range("A27").Style = "correctStyle"
$ws.Range("A28").Value = "Option D"
$ws.Range("A28").Style = "correctStyle"
$ws.Range("A29").Value = "Option B"
$ws.Range("A29").Style = "incorrectStyle"
$ws.Range("A31").Value = "Option B"
$ws.Range("A31").Style = "incorrectStyle"
$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"
$ws.Range("A33").Value = "Option D"
$ws.Range("A33").Style = "correctStyle"
$ws.Range("A34").Value = "Option B"
$ws.Range("A34").Style = "correctStyle"
$ws.Range("A35").Value = "Option B"
$ws.Range("A35").Style = "incorrectStyle"
$ws.Range("A36").Value = "Option A"
$ws.Range("A36").Style = "correctStyle"
$ws.Range("A37").Value = "Option B"
$ws.Range("A37").Style = "incorrectStyle"
$ws.Range("A38").Value = "Option B"
$ws.Range("A38").Style = "incorrectStyle"
$ws.Range("A39").Value = "Option D"
$ws.Range("A39").Style = "correctStyle"
$ws.Range("A40").Value = "Option D"
$ws.Range("A40").Style = "correctStyle"
